# Update developer guide images#3
#
# 1. Refresh the cached "Update automatically" date placeholder text
#    (10/16/2016 -> 3/17/2018) on every slide layout, the slide master,
#    and the notes master.
# 2. Rename a few method/event identifiers used in the sequence-diagram
#    textboxes on slide 1.

$p = $ppt.ActivePresentation

$oldDate = "10/16/2016"
$newDate = "3/17/2018"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# -- Slide layouts (11 of them) --
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# -- Slide master --
Update-DatePlaceholder $p.SlideMaster.Shapes

# -- Notes master --
Update-DatePlaceholder $p.NotesMaster.Shapes

# -- Slide 1 sequence-diagram textbox renames --
$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if (-not $sh.HasTextFrame) { continue }
    $tr = $sh.TextFrame.TextRange
    $txt = $tr.Text

    if ($txt -eq "deletePerson(p)") {
        $tr.Text = "deleteMember(p)"
    }
    elseif ($txt -eq "post(AddressBookChangedEvent)") {
        $tr.Text = "post(ClubBookChangedEvent)"
    }
    elseif ($txt -eq "handleAddresssBookChangedEvent()") {
        $tr.Text = "handleClubBookChangedEvent()"
    }
}
